$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.576682333333332
$ws.Range("H2").Value = 25.730047
$ws.Range("I2").Value = 0.3754230651280642
$ws.Range("J2").Value = 0.3754230651280643
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.39091433333333
$ws.Range("N2").Value = 31.172743
$ws.Range("O2").Value = 0.0835098648954196
$ws.Range("P2").Value = 0.0835098648954196
$ws.Range("Q2").Value = 89.11957138988009
$ws.Range("R2").Value = 802.0761425089208
$ws.Range("S2").Value = 0.03135152944746896
$ws.Range("T2").Value = 0.03135152944746896
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.576682333333332
$ws.Range("H3").Value = 25.730047
$ws.Range("I3").Value = 0.3754230651280642
$ws.Range("J3").Value = 0.3754230651280643
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.71561933333334
$ws.Range("N3").Value = 107.146858
$ws.Range("O3").Value = 0.287039855156433
$ws.Range("P3").Value = 0.287039855156433
$ws.Range("Q3").Value = 306.3215213602584
$ws.Range("R3").Value = 2756.893692242326
$ws.Range("S3").Value = 0.1077613822367437
$ws.Range("T3").Value = 0.1077613822367437
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.576682333333332
$ws.Range("H4").Value = 25.730047
$ws.Range("I4").Value = 0.3754230651280642
$ws.Range("J4").Value = 0.3754230651280643
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 51.87044599999999
$ws.Range("N4").Value = 155.611338
$ws.Range("O4").Value = 0.4168732219867682
$ws.Range("P4").Value = 0.4168732219867682
$ws.Range("Q4").Value = 444.8763378303206
$ws.Range("R4").Value = 4003.887040472885
$ws.Range("S4").Value = 0.1565038227680844
$ws.Range("T4").Value = 0.1565038227680845
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.576682333333332
$ws.Range("H5").Value = 25.730047
$ws.Range("I5").Value = 0.3754230651280642
$ws.Range("J5").Value = 0.3754230651280643
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 26.450408
$ws.Range("N5").Value = 79.351224
$ws.Range("O5").Value = 0.2125770579613792
$ws.Range("P5").Value = 0.2125770579613792
$ws.Range("Q5").Value = 226.8567470030586
$ws.Range("R5").Value = 2041.710723027528
$ws.Range("S5").Value = 0.07980633067576715
$ws.Range("T5").Value = 0.07980633067576717
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.028452333333334
$ws.Range("H6").Value = 21.085357
$ws.Range("I6").Value = 0.3076531245457689
$ws.Range("J6").Value = 0.3076531245457689
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.39091433333333
$ws.Range("N6").Value = 31.172743
$ws.Range("O6").Value = 0.0835098648954196
$ws.Range("P6").Value = 0.0835098648954196
$ws.Range("Q6").Value = 73.03204609158345
$ws.Range("R6").Value = 657.288414824251
$ws.Range("S6").Value = 0.02569207086547086
$ws.Range("T6").Value = 0.02569207086547086
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.028452333333334
$ws.Range("H7").Value = 21.085357
$ws.Range("I7").Value = 0.3076531245457689
$ws.Range("J7").Value = 0.3076531245457689
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 35.71561933333334
$ws.Range("N7").Value = 107.146858
$ws.Range("O7").Value = 0.287039855156433
$ws.Range("P7").Value = 0.287039855156433
$ws.Range("Q7").Value = 251.0255280398118
$ws.Range("R7").Value = 2259.229752358307
$ws.Range("S7").Value = 0.08830870830804154
$ws.Range("T7").Value = 0.08830870830804154
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.028452333333334
$ws.Range("H8").Value = 21.085357
$ws.Range("I8").Value = 0.3076531245457689
$ws.Range("J8").Value = 0.3076531245457689
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 51.87044599999999
$ws.Range("N8").Value = 155.611338
$ws.Range("O8").Value = 0.4168732219867682
$ws.Range("P8").Value = 0.4168732219867682
$ws.Range("Q8").Value = 364.5689572197406
$ws.Range("R8").Value = 3281.120614977666
$ws.Range("S8").Value = 0.1282523492836911
$ws.Range("T8").Value = 0.1282523492836912
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.028452333333334
$ws.Range("H9").Value = 21.085357
$ws.Range("I9").Value = 0.3076531245457689
$ws.Range("J9").Value = 0.3076531245457689
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.450408
$ws.Range("N9").Value = 79.351224
$ws.Range("O9").Value = 0.2125770579613792
$ws.Range("P9").Value = 0.2125770579613792
$ws.Range("Q9").Value = 185.9054318252187
$ws.Range("R9").Value = 1673.148886426968
$ws.Range("S9").Value = 0.06539999608856532
$ws.Range("T9").Value = 0.06539999608856534
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4119716666666666
$ws.Range("H10").Value = 1.235915
$ws.Range("I10").Value = 0.01803304119645609
$ws.Range("J10").Value = 0.01803304119645609
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.39091433333333
$ws.Range("N10").Value = 31.172743
$ws.Range("O10").Value = 0.0835098648954196
$ws.Range("P10").Value = 0.0835098648954196
$ws.Range("Q10").Value = 4.280762296093888
$ws.Range("R10").Value = 38.52686066484499
$ws.Range("S10").Value = 0.001505936833969584
$ws.Range("T10").Value = 0.001505936833969584
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4119716666666666
$ws.Range("H11").Value = 1.235915
$ws.Range("I11").Value = 0.01803304119645609
$ws.Range("J11").Value = 0.01803304119645609
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 35.71561933333334
$ws.Range("N11").Value = 107.146858
$ws.Range("O11").Value = 0.287039855156433
$ws.Range("P11").Value = 0.287039855156433
$ws.Range("Q11").Value = 14.71382322278555
$ws.Range("R11").Value = 132.42440900507
$ws.Range("S11").Value = 0.005176201533060745
$ws.Range("T11").Value = 0.005176201533060747
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4119716666666666
$ws.Range("H12").Value = 1.235915
$ws.Range("I12").Value = 0.01803304119645609
$ws.Range("J12").Value = 0.01803304119645609
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 51.87044599999999
$ws.Range("N12").Value = 155.611338
$ws.Range("O12").Value = 0.4168732219867682
$ws.Range("P12").Value = 0.4168732219867682
$ws.Range("Q12").Value = 21.36915408936333
$ws.Range("R12").Value = 192.32238680427
$ws.Range("S12").Value = 0.007517491985786775
$ws.Range("T12").Value = 0.007517491985786778
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4119716666666666
$ws.Range("H13").Value = 1.235915
$ws.Range("I13").Value = 0.01803304119645609
$ws.Range("J13").Value = 0.01803304119645609
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 26.450408
$ws.Range("N13").Value = 79.351224
$ws.Range("O13").Value = 0.2125770579613792
$ws.Range("P13").Value = 0.2125770579613792
$ws.Range("Q13").Value = 10.89681866777333
$ws.Range("R13").Value = 98.07136800996
$ws.Range("S13").Value = 0.003833410843638986
$ws.Range("T13").Value = 0.003833410843638987
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.828273
$ws.Range("H14").Value = 20.484819
$ws.Range("I14").Value = 0.2988907691297108
$ws.Range("J14").Value = 0.2988907691297109
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.39091433333333
$ws.Range("N14").Value = 31.172743
$ws.Range("O14").Value = 0.0835098648954196
$ws.Range("P14").Value = 0.0835098648954196
$ws.Range("Q14").Value = 70.951999787613
$ws.Range("R14").Value = 638.567998088517
$ws.Range("S14").Value = 0.0249603277485102
$ws.Range("T14").Value = 0.02496032774851021
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.828273
$ws.Range("H15").Value = 20.484819
$ws.Range("I15").Value = 0.2988907691297108
$ws.Range("J15").Value = 0.2988907691297109
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 35.71561933333334
$ws.Range("N15").Value = 107.146858
$ws.Range("O15").Value = 0.287039855156433
$ws.Range("P15").Value = 0.287039855156433
$ws.Range("Q15").Value = 243.875999172078
$ws.Range("R15").Value = 2194.883992548702
$ws.Range("S15").Value = 0.08579356307858704
$ws.Range("T15").Value = 0.08579356307858706
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.828273
$ws.Range("H16").Value = 20.484819
$ws.Range("I16").Value = 0.2988907691297108
$ws.Range("J16").Value = 0.2988907691297109
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 51.87044599999999
$ws.Range("N16").Value = 155.611338
$ws.Range("O16").Value = 0.4168732219867682
$ws.Range("P16").Value = 0.4168732219867682
$ws.Range("Q16").Value = 354.1855659197579
$ws.Range("R16").Value = 3187.670093277822
$ws.Range("S16").Value = 0.1245995579492058
$ws.Range("T16").Value = 0.1245995579492059
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.828273
$ws.Range("H17").Value = 20.484819
$ws.Range("I17").Value = 0.2988907691297108
$ws.Range("J17").Value = 0.2988907691297109
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 26.450408
$ws.Range("N17").Value = 79.351224
$ws.Range("O17").Value = 0.2125770579613792
$ws.Range("P17").Value = 0.2125770579613792
$ws.Range("Q17").Value = 180.610606785384
$ws.Range("R17").Value = 1625.495461068456
$ws.Range("S17").Value = 0.06353732035340774
$ws.Range("T17").Value = 0.06353732035340777
